$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "30.935.23"
Set-TextValue $ws.Range("E2") "  +3.69%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.679.06"
Set-TextValue $ws.Range("E3") "  +3.22%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.15%  "

# Row 5
Set-TextValue $ws.Range("D5") "219.59"

# Row 6
Set-TextValue $ws.Range("D6") "0.533"
Set-TextValue $ws.Range("E6") "  +2.81%  "

# Row 7
Set-TextValue $ws.Range("E7") "  -0.14%  "

# Row 8
Set-TextValue $ws.Range("D8") "29.08"
Set-TextValue $ws.Range("E8") "  +2.12%  "

# Row 9
Set-TextValue $ws.Range("E9") "  +2.76%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0641"
Set-TextValue $ws.Range("E10") "  +5.52%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0909"
Set-TextValue $ws.Range("E11") "  +1.08%  "

# Row 12
Set-TextValue $ws.Range("D12") "1.920.42"
Set-TextValue $ws.Range("E12") "  +3.27%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.687.01"
Set-TextValue $ws.Range("E13") "  +3.49%  "

# Row 14
Set-TextValue $ws.Range("D14") "0.604"
Set-TextValue $ws.Range("E14") "  +7.15%  "

# Row 15
Set-TextValue $ws.Range("D15") "10.05"
Set-TextValue $ws.Range("E15") "  +8.42%  "

# Row 16
Set-TextValue $ws.Range("D16") "4.10"
Set-TextValue $ws.Range("E16") "  +7.18%  "

# Row 17
Set-TextValue $ws.Range("D17") "30.896.93"
Set-TextValue $ws.Range("E17") "  +3.47%  "

# Row 18
Set-TextValue $ws.Range("D18") "66.05"
Set-TextValue $ws.Range("E18") "  +2.17%  "

# Row 19
Set-TextValue $ws.Range("D19") "246.33"
Set-TextValue $ws.Range("E19") "  +2.61%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.0₃0719"
Set-TextValue $ws.Range("E20") "  +2.68%  "

# Row 21
Set-TextValue $ws.Range("E21") "  -0.05%  "

# Row 22
Set-TextValue $ws.Range("E22") "  +3.36%  "

# Row 23
Set-TextValue $ws.Range("D23") "9.96"
Set-TextValue $ws.Range("E23") "  +2.04%  "

# Row 24
Set-TextValue $ws.Range("E24") "  -0.22%  "

# Row 25
Set-TextValue $ws.Range("D25") "159.22"
Set-TextValue $ws.Range("E25") "  +1.05%  "

# Row 26
Set-TextValue $ws.Range("E26") "  +2.50%  "

# Row 27
Set-TextValue $ws.Range("E27") "  +2.46%  "

# Row 28
Set-TextValue $ws.Range("E28") "  +2.07%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.00"
Set-TextValue $ws.Range("E29") "  -0.01%  "

# Row 30
Set-TextValue $ws.Range("E30") "  +0.98%  "

# Row 31
Set-TextValue $ws.Range("E31") "  +3.58%  "

# Row 32
Set-TextValue $ws.Range("D32") "3.47"
Set-TextValue $ws.Range("E32") "  +3.36%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.519.20"
Set-TextValue $ws.Range("E33") "  +6.59%  "

# Row 34
Set-TextValue $ws.Range("D34") "3.32"
Set-TextValue $ws.Range("E34") "  +4.99%  "

# Row 35
Set-TextValue $ws.Range("E35") "  +4.48%  "

# Row 36
Set-TextValue $ws.Range("D36") "84.53"
Set-TextValue $ws.Range("E36") "  +12.87%  "

# Row 37
Set-TextValue $ws.Range("E37") "  +0.73%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.604"
Set-TextValue $ws.Range("E38") "  +8.84%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.0180"
Set-TextValue $ws.Range("E39") "  +5.51%  "

# Row 40
Set-TextValue $ws.Range("B40") "HuobiToken"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D40") "2.29"
Set-TextValue $ws.Range("E40") "  +0.15%  "

# Row 41
Set-TextValue $ws.Range("B41") "MXToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D41") "2.64"
Set-TextValue $ws.Range("E41") "  -3.89%  "

# Row 42
Set-TextValue $ws.Range("D42") "2.05"
Set-TextValue $ws.Range("E42") "  +3.54%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.838"
Set-TextValue $ws.Range("E43") "  +1.28%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.0502"
Set-TextValue $ws.Range("E44") "  +0.86%  "

# Row 45
Set-TextValue $ws.Range("E45") "  +2.18%  "

# Row 46
Set-TextValue $ws.Range("E46") "  +0.00%  "

# Row 47
Set-TextValue $ws.Range("E47") "  +4.80%  "

# Row 48
Set-TextValue $ws.Range("D48") "51.02"

# Row 49
Set-TextValue $ws.Range("D49") "1.814.28"
Set-TextValue $ws.Range("E49") "  +2.70%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.0₆0118"
Set-TextValue $ws.Range("E50") "  +6.78%  "

# Row 51
Set-TextValue $ws.Range("D51") "93.08"
Set-TextValue $ws.Range("E51") "  +2.59%  "
